# Refresh the cryptos list: updates Price (col D) and Volume(1h) (col E) for
# most rows, and for rows 43-44 also swaps Hedera/dogwifhat (Coin/Link/Price/
# Volume) with their refreshed figures.
#
# Price cells are stored as plain text (e.g. "581.96", "67.941.10") even when
# they look numeric, so for any new price that parses as a float we force the
# cell to Text format before writing, then restore the default "Normal" style
# so no stray number-format style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.941.10'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '3.261.22'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').Value = '3.260.48'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('E10').Value = '  -3.01%  '
$ws.Range('E11').Value = '  -2.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.412'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '3.826.85'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('D16').Value = '67.850.83'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '3.260.11'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '394.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.516'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000118'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('E27').Value = '  -2.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.14%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.54'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.74%  '
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.94%  '
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('E40').Value = '  -3.43%  '
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.45%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0689'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').Value = '2.613.72'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '334.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0277'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('E51').Value = '  -0.86%  '
